# Apply corrections to the SoIB_summaries.xlsx workbook after fixing
# error estimation and the number of projected years.

$wb = $excel.ActiveWorkbook

# --- Sheet: "Trends Status" ---
$wsTrends = $wb.Worksheets.Item("Trends Status")

$wsTrends.Range("C2").Value = 2

$wsTrends.Range("C3").Value = 2
$wsTrends.Range("D3").Value = 25

$wsTrends.Range("B6").Value = 3
$wsTrends.Range("D6").Value = 75

$wsTrends.Range("B7").Value = 21
$wsTrends.Range("C7").Value = 26

# --- Sheet: "Species qualification" ---
$wsSpecies = $wb.Worksheets.Item("Species qualification")

$wsSpecies.Range("C3").Value = 4
$wsSpecies.Range("C4").Value = 4
